$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I0 (col I) and IF (col J) in row 1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy header formatting (bold, border, centered) from H1 onto the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data for column I (rows 2-73) and column J (rows 2-73)
$iValues = @(5,4,5,6,7,7,9,7,8,6,6,6,8,10,8,10,5,8,8,8,5,7,9,6,6,9,7,5,6,6,5,8,5,6,7,7,5,7,7,9,6,2,6,5,8,9,10,5,8,6,7,8,4,6,10,8,7,6,8,8,9,7,6,6,5,5,4,5,4,5,4,3)
$jValues = @(5,5,5,6,7,8,9,7,9,6,7,6,8,10,8,10,6,8,8,8,5,7,9,6,6,9,7,5,6,7,6,8,6,6,7,7,6,7,8,9,6,3,6,6,8,9,10,5,8,6,7,8,5,6,10,8,7,6,8,8,9,7,7,7,5,5,4,5,5,5,4,3)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}

Write-Host "Applied I0/IF columns"
